$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.6879308622404
$ws.Range("C2").Value = 0.254893879407089
$ws.Range("B3").Value = 0.130167622053188
$ws.Range("C3").Value = 0.11158061989941

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.91397156346093
$ws.Range("C2").Value = 0.251259573281716
$ws.Range("B3").Value = -0.965044938060798
$ws.Range("C3").Value = 0.0796827302863566

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.0339517655342
$ws.Range("C2").Value = 0.179930231478184
$ws.Range("B3").Value = 1.64617655510194
$ws.Range("C3").Value = 0.157352972748654

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.58753790062641
$ws.Range("C2").Value = 0.210820747593303
$ws.Range("B3").Value = 0.0146205057425297
$ws.Range("C3").Value = 0.0191134696191717

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0649708897591956
$ws.Range("B2").Value = -0.0206858569923194
$ws.Range("A3").Value = -0.0206858569923194
$ws.Range("B3").Value = 0.0124502347371367

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0631313731657098
$ws.Range("B2").Value = -0.0158753027847061
$ws.Range("A3").Value = -0.0158753027847061
$ws.Range("B3").Value = 0.00634933750588825

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0323748881997928
$ws.Range("B2").Value = -0.00937185231836174
$ws.Range("A3").Value = -0.00937185231836174
$ws.Range("B3").Value = 0.0247599580328386

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0444453876157992
$ws.Range("B2").Value = -0.00225609370786197
$ws.Range("A3").Value = -0.00225609370786197
$ws.Range("B3").Value = 0.000365324720882999
